$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "f_count" (row 11) and "p_count" (row 12) rows entirely.
#    This shifts the former rows 13-19 (toa, distance, c_walls, w_walls,
#    exp_pl, n_power, esp) up to become rows 11-17, carrying their
#    numeric statistics (C:I) along unchanged.
$ws.Range("A11:A12").EntireRow.Delete()

# 2. Add the new "Field" header cell in A1, cloning the header style (bold,
#    bordered, centered) already used by the other row-1 header cells.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Rewrite the header row text.
$ws.Range("A1").Value2 = "Field"
$ws.Range("B1").Value2 = "Unit"
$ws.Range("C1").Value2 = "Mean"
$ws.Range("D1").Value2 = "STD (" + [string][char]0x03C3 + ")"
$ws.Range("E1").Value2 = "Min"
$ws.Range("F1").Value2 = "1" + [string][char]0x02E2 + [string][char]0x1D57 + " Quartile (25%)"
$ws.Range("G1").Value2 = "Median"
$ws.Range("H1").Value2 = "3" + [string][char]0x02B3 + [string][char]0x1D48 + " Quartile (75%)"
$ws.Range("I1").Value2 = "Max"

# 4. Column B no longer holds the observation "count" - it now holds the
#    unit of each field. Replace the numeric counts with unit strings.
$units = @{
    2  = "ppm"                                                     # co2
    3  = "%"                                                       # humidity
    4  = [string][char]0x00B5 + "g/m" + [string][char]0x00B3       # pm25 -> ug/m3
    5  = "hPa"                                                     # pressure
    6  = [string][char]0x00B0 + "C"                                # temperature -> degC
    7  = "dBm"                                                     # rssi
    8  = "dB"                                                      # snr
    9  = "bit/sym"                                                 # SF
    10 = "MHz"                                                     # frequency
    11 = "s"                                                       # toa
    12 = "m"                                                       # distance
    13 = ""                                                        # c_walls (no unit)
    14 = ""                                                        # w_walls (no unit)
    15 = "dB"                                                      # exp_pl
    16 = "dB"                                                      # n_power
    17 = "dBm"                                                     # esp
}
foreach ($row in $units.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $units[$row]
}

# 5. The field-name cells in column A (rows 2-17) are no longer bold /
#    bordered - only the header row keeps that styling now.
$ws.Range("A2:A17").ClearFormats()
